# [UPDATE] last details footer and navBar
$wb = $excel.ActiveWorkbook

# --- "Personas" sheet: stretch the used range down to row 219 (30 extra
#     padding rows in the trailing marker column, mirroring the existing
#     Q3:Q189 footer/navBar padding pattern) ---
$wsPersonas = $wb.Worksheets.Item("Personas")
$wsPersonas.Range("Q190:Q219").Value2 = 0

# --- "Subcategorias" sheet: the old "Becados" sub-category is replaced by
#     "Becados LF" (and its description/owner), the stray test rows below
#     it are removed ---
$wsSub = $wb.Worksheets.Item("Subcategorias")

# Drop rows 4:7 (the Ingles B2 / Mantenimientos / Emergencias / test rows)
$wsSub.Range("A4:A7").EntireRow.Delete()

# Refresh row 3 in place
$wsSub.Range("A3").Value2 = 59
$wsSub.Range("B3").Value2 = "Becados LF"
$wsSub.Range("C3").Value2 = "Becados LF Monto económico"
$wsSub.Range("D3").Value2 = "Herminia Ávila"

$wsSub.Range("A3:D3").Select()
